$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
# F2: 1191 -> 1195
$ws.Range("F2").Value = 1195
# F4: 48 -> 49
$ws.Range("F4").Value = 49
# F5: 1312 -> 1329
$ws.Range("F5").Value = 1329
# F6: 1696 -> 1701
$ws.Range("F6").Value = 1701
# F7: 6214 -> 6217
$ws.Range("F7").Value = 6217
# F8: 124 -> 125
$ws.Range("F8").Value = 125
# F9: 1814 -> 1818
$ws.Range("F9").Value = 1818
# F12: 15 -> 16
$ws.Range("F12").Value = 16
# F15: 21 -> 22
$ws.Range("F15").Value = 22
# F16: 6908 -> 6925
$ws.Range("F16").Value = 6925
# F17: 124 -> 125
$ws.Range("F17").Value = 125
# F18: 53 -> 54
$ws.Range("F18").Value = 54
# F19: 164 -> 165
$ws.Range("F19").Value = 165
# F21: 1703 -> 1704
$ws.Range("F21").Value = 1704
# F23: 14 -> 15
$ws.Range("F23").Value = 15
# F26: 1576 -> 1580
$ws.Range("F26").Value = 1580
# F27: 756 -> 757
$ws.Range("F27").Value = 757
# F28: 313 -> 315
$ws.Range("F28").Value = 315
# F31: 52 -> 54
$ws.Range("F31").Value = 54
# F33: 3891 -> 3893
$ws.Range("F33").Value = 3893

$ws = $wb.Worksheets.Item("演出")
# F24: 1 -> 2
$ws.Range("F24").Value = 2

$ws = $wb.Worksheets.Item("本地生活")
# F2: 9527 -> 9529
$ws.Range("F2").Value = 9529
# F5: 242 -> 243
$ws.Range("F5").Value = 243

$ws = $wb.Worksheets.Item("全部类型")
# F2: 9527 -> 9529
$ws.Range("F2").Value = 9529
# F5: 1191 -> 1195
$ws.Range("F5").Value = 1195
# F7: 48 -> 49
$ws.Range("F7").Value = 49
# F10: 1312 -> 1329
$ws.Range("F10").Value = 1329
# F11: 242 -> 243
$ws.Range("F11").Value = 243
# F12: 1696 -> 1701
$ws.Range("F12").Value = 1701
# F13: 6214 -> 6217
$ws.Range("F13").Value = 6217
# F14: 124 -> 125
$ws.Range("F14").Value = 125
# F15: 1814 -> 1818
$ws.Range("F15").Value = 1818
# F20: 15 -> 16
$ws.Range("F20").Value = 16
# F23: 6908 -> 6925
$ws.Range("F23").Value = 6925
# F24: 124 -> 125
$ws.Range("F24").Value = 125
# F25: 53 -> 54
$ws.Range("F25").Value = 54
# F26: 164 -> 165
$ws.Range("F26").Value = 165
# F28: 1703 -> 1704
$ws.Range("F28").Value = 1704
# F30: 14 -> 15
$ws.Range("F30").Value = 15
# F33: 1576 -> 1580
$ws.Range("F33").Value = 1580
# F34: 756 -> 757
$ws.Range("F34").Value = 757
# F36: 313 -> 315
$ws.Range("F36").Value = 315
# F45: 3891 -> 3893
$ws.Range("F45").Value = 3893
